$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily data rows appended for 15, 16, 17 April 2021 (serials 44301-44303)
$data = @(
    @{ Row = 227; A = 44301; B = 0; C = 7;  D = 113.2136503315543 },
    @{ Row = 228; A = 44302; B = 6; C = 12; D = 194.0805434255216 },
    @{ Row = 229; A = 44303; B = 2; C = 13; D = 210.2539220443151 }
)

foreach ($entry in $data) {
    $r = $entry.Row

    # Write the values first
    $ws.Cells.Item($r, 1).Value = $entry.A
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D

    # Copy the date-column formatting (style index 2: centered, bordered,
    # date/time number format) from the row above, like the existing rows.
    $ws.Cells.Item($r - 1, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
}

$excel.CutCopyMode = $false
